$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -17.33501518566552
$ws.Cells.Item(2, 3).Value = 4.321919805167666
$ws.Cells.Item(2, 4).Value = -17.33501518566552
$ws.Cells.Item(2, 5).Value = -17.33501518566552
$ws.Cells.Item(2, 6).Value = -17.33501518566552
$ws.Cells.Item(2, 7).Value = -17.33501518566552
$ws.Cells.Item(2, 8).Value = -17.33501518566552
$ws.Cells.Item(2, 9).Value = -17.33501518566552
$ws.Cells.Item(2, 10).Value = -17.33501518566552
$ws.Cells.Item(2, 11).Value = -17.33501518566552
$ws.Cells.Item(3, 2).Value = -17.33501518566552
$ws.Cells.Item(3, 3).Value = -17.33501518566552
$ws.Cells.Item(3, 4).Value = -17.33501518566552
$ws.Cells.Item(3, 5).Value = -17.33501518566552
$ws.Cells.Item(3, 6).Value = -17.33501518566552
$ws.Cells.Item(3, 7).Value = -17.33501518566552
$ws.Cells.Item(3, 8).Value = -17.33501518566552
$ws.Cells.Item(3, 9).Value = 2.688789242297438
$ws.Cells.Item(3, 10).Value = -17.33501518566552
$ws.Cells.Item(3, 11).Value = -17.33501518566552
$ws.Cells.Item(4, 2).Value = -17.33501518566552
$ws.Cells.Item(4, 3).Value = -17.33501518566552
$ws.Cells.Item(4, 4).Value = 1.638471989017111
$ws.Cells.Item(4, 5).Value = -17.33501518566552
$ws.Cells.Item(4, 6).Value = 3.415846335912579
$ws.Cells.Item(4, 7).Value = -17.33501518566552
$ws.Cells.Item(4, 8).Value = 1.544200092457836
$ws.Cells.Item(4, 9).Value = -17.33501518566552
$ws.Cells.Item(4, 10).Value = 2.238545099672264
$ws.Cells.Item(4, 11).Value = -17.33501518566552
$ws.Cells.Item(5, 2).Value = -17.33501518566552
$ws.Cells.Item(5, 3).Value = -17.33501518566552
$ws.Cells.Item(5, 4).Value = -17.33501518566552
$ws.Cells.Item(5, 5).Value = -17.33501518566552
$ws.Cells.Item(5, 6).Value = -17.33501518566552
$ws.Cells.Item(5, 7).Value = 2.869007823297684
$ws.Cells.Item(5, 8).Value = -17.33501518566552
$ws.Cells.Item(5, 9).Value = -17.33501518566552
$ws.Cells.Item(5, 10).Value = -17.33501518566552
$ws.Cells.Item(5, 11).Value = -17.33501518566552
$ws.Cells.Item(6, 2).Value = -17.33501518566552
$ws.Cells.Item(6, 3).Value = -17.33501518566552
$ws.Cells.Item(6, 4).Value = -17.33501518566552
$ws.Cells.Item(6, 5).Value = -17.33501518566552
$ws.Cells.Item(6, 6).Value = -17.33501518566552
$ws.Cells.Item(6, 7).Value = -17.33501518566552
$ws.Cells.Item(6, 8).Value = -17.33501518566552
$ws.Cells.Item(6, 9).Value = -17.33501518566552
$ws.Cells.Item(6, 10).Value = -17.33501518566552
$ws.Cells.Item(6, 11).Value = -17.33501518566552
$ws.Cells.Item(7, 2).Value = 2.814008708460704
$ws.Cells.Item(7, 3).Value = -17.33501518566552
$ws.Cells.Item(7, 4).Value = -17.33501518566552
$ws.Cells.Item(7, 5).Value = -17.33501518566552
$ws.Cells.Item(7, 6).Value = -17.33501518566552
$ws.Cells.Item(7, 7).Value = -17.33501518566552
$ws.Cells.Item(7, 8).Value = -17.33501518566552
$ws.Cells.Item(7, 9).Value = -17.33501518566552
$ws.Cells.Item(7, 10).Value = -17.33501518566552
$ws.Cells.Item(7, 11).Value = -17.33501518566552
$ws.Cells.Item(8, 2).Value = -17.33501518566552
$ws.Cells.Item(8, 3).Value = -17.33501518566552
$ws.Cells.Item(8, 4).Value = -17.33501518566552
$ws.Cells.Item(8, 5).Value = 1.312362429840031
$ws.Cells.Item(8, 6).Value = -17.33501518566552
$ws.Cells.Item(8, 7).Value = -17.33501518566552
$ws.Cells.Item(8, 8).Value = -17.33501518566552
$ws.Cells.Item(8, 9).Value = -17.33501518566552
$ws.Cells.Item(8, 10).Value = -17.33501518566552
$ws.Cells.Item(8, 11).Value = -17.33501518566552
$ws.Cells.Item(9, 2).Value = 3.696832046162206
$ws.Cells.Item(9, 3).Value = -17.33501518566552
$ws.Cells.Item(9, 4).Value = -17.33501518566552
$ws.Cells.Item(9, 5).Value = -17.33501518566552
$ws.Cells.Item(9, 6).Value = -17.33501518566552
$ws.Cells.Item(9, 7).Value = -17.33501518566552
$ws.Cells.Item(9, 8).Value = -17.33501518566552
$ws.Cells.Item(9, 9).Value = -17.33501518566552
$ws.Cells.Item(9, 10).Value = -17.33501518566552
$ws.Cells.Item(9, 11).Value = -17.33501518566552
$ws.Cells.Item(10, 2).Value = -17.33501518566552
$ws.Cells.Item(10, 3).Value = -17.33501518566552
$ws.Cells.Item(10, 4).Value = -17.33501518566552
$ws.Cells.Item(10, 5).Value = -17.33501518566552
$ws.Cells.Item(10, 6).Value = -17.33501518566552
$ws.Cells.Item(10, 7).Value = -17.33501518566552
$ws.Cells.Item(10, 8).Value = -17.33501518566552
$ws.Cells.Item(10, 9).Value = 1.273618657876925
$ws.Cells.Item(10, 10).Value = -17.33501518566552
$ws.Cells.Item(10, 11).Value = 1.872087335926466
$ws.Cells.Item(11, 2).Value = -17.33501518566552
$ws.Cells.Item(11, 3).Value = -17.33501518566552
$ws.Cells.Item(11, 4).Value = -17.33501518566552
$ws.Cells.Item(11, 5).Value = 2.652694325184395
$ws.Cells.Item(11, 6).Value = -17.33501518566552
$ws.Cells.Item(11, 7).Value = 2.773183813319917
$ws.Cells.Item(11, 8).Value = -17.33501518566552
$ws.Cells.Item(11, 9).Value = -17.33501518566552
$ws.Cells.Item(11, 10).Value = -17.33501518566552
$ws.Cells.Item(11, 11).Value = 1.752667776584454
$ws.Cells.Item(12, 2).Value = -17.33501518566552
$ws.Cells.Item(12, 3).Value = -17.33501518566552
$ws.Cells.Item(12, 4).Value = -17.33501518566552
$ws.Cells.Item(12, 5).Value = -17.33501518566552
$ws.Cells.Item(12, 6).Value = -17.33501518566552
$ws.Cells.Item(12, 7).Value = -17.33501518566552
$ws.Cells.Item(12, 8).Value = -17.33501518566552
$ws.Cells.Item(12, 9).Value = -17.33501518566552
$ws.Cells.Item(12, 10).Value = -17.33501518566552
$ws.Cells.Item(12, 11).Value = -17.33501518566552
$ws.Cells.Item(13, 2).Value = -17.33501518566552
$ws.Cells.Item(13, 3).Value = -17.33501518566552
$ws.Cells.Item(13, 4).Value = -17.33501518566552
$ws.Cells.Item(13, 5).Value = 2.421282139949281
$ws.Cells.Item(13, 6).Value = -17.33501518566552
$ws.Cells.Item(13, 7).Value = -17.33501518566552
$ws.Cells.Item(13, 8).Value = -17.33501518566552
$ws.Cells.Item(13, 9).Value = -17.33501518566552
$ws.Cells.Item(13, 10).Value = 2.030369604175275
$ws.Cells.Item(13, 11).Value = 1.853518124640224
$ws.Cells.Item(14, 2).Value = -17.33501518566552
$ws.Cells.Item(14, 3).Value = -17.33501518566552
$ws.Cells.Item(14, 4).Value = 1.428758030840785
$ws.Cells.Item(14, 5).Value = -17.33501518566552
$ws.Cells.Item(14, 6).Value = -17.33501518566552
$ws.Cells.Item(14, 7).Value = -17.33501518566552
$ws.Cells.Item(14, 8).Value = -17.33501518566552
$ws.Cells.Item(14, 9).Value = -17.33501518566552
$ws.Cells.Item(14, 10).Value = -17.33501518566552
$ws.Cells.Item(14, 11).Value = 2.152253725898614
$ws.Cells.Item(15, 2).Value = -17.33501518566552
$ws.Cells.Item(15, 3).Value = -17.33501518566552
$ws.Cells.Item(15, 4).Value = 1.23031378388058
$ws.Cells.Item(15, 5).Value = -17.33501518566552
$ws.Cells.Item(15, 6).Value = -17.33501518566552
$ws.Cells.Item(15, 7).Value = -17.33501518566552
$ws.Cells.Item(15, 8).Value = -17.33501518566552
$ws.Cells.Item(15, 9).Value = -17.33501518566552
$ws.Cells.Item(15, 10).Value = -17.33501518566552
$ws.Cells.Item(15, 11).Value = -17.33501518566552
$ws.Cells.Item(16, 2).Value = -17.33501518566552
$ws.Cells.Item(16, 3).Value = -17.33501518566552
$ws.Cells.Item(16, 4).Value = -17.33501518566552
$ws.Cells.Item(16, 5).Value = -17.33501518566552
$ws.Cells.Item(16, 6).Value = -17.33501518566552
$ws.Cells.Item(16, 7).Value = -17.33501518566552
$ws.Cells.Item(16, 8).Value = -17.33501518566552
$ws.Cells.Item(16, 9).Value = -17.33501518566552
$ws.Cells.Item(16, 10).Value = 2.188081534657647
$ws.Cells.Item(16, 11).Value = -17.33501518566552
$ws.Cells.Item(17, 2).Value = -17.33501518566552
$ws.Cells.Item(17, 3).Value = -17.33501518566552
$ws.Cells.Item(17, 4).Value = 2.042282205502006
$ws.Cells.Item(17, 5).Value = -17.33501518566552
$ws.Cells.Item(17, 6).Value = -17.33501518566552
$ws.Cells.Item(17, 7).Value = -17.33501518566552
$ws.Cells.Item(17, 8).Value = 1.223753034857074
$ws.Cells.Item(17, 9).Value = 1.861272354323441
$ws.Cells.Item(17, 10).Value = 1.922724834507435
$ws.Cells.Item(17, 11).Value = -17.33501518566552
$ws.Cells.Item(18, 2).Value = -17.33501518566552
$ws.Cells.Item(18, 3).Value = -17.33501518566552
$ws.Cells.Item(18, 4).Value = -17.33501518566552
$ws.Cells.Item(18, 5).Value = -17.33501518566552
$ws.Cells.Item(18, 6).Value = -17.33501518566552
$ws.Cells.Item(18, 7).Value = -17.33501518566552
$ws.Cells.Item(18, 8).Value = 1.472933986030219
$ws.Cells.Item(18, 9).Value = 1.159891168450407
$ws.Cells.Item(18, 10).Value = 1.509558801285808
$ws.Cells.Item(18, 11).Value = -17.33501518566552
$ws.Cells.Item(19, 2).Value = -17.33501518566552
$ws.Cells.Item(19, 3).Value = -17.33501518566552
$ws.Cells.Item(19, 4).Value = 1.879552958418121
$ws.Cells.Item(19, 5).Value = -17.33501518566552
$ws.Cells.Item(19, 6).Value = -17.33501518566552
$ws.Cells.Item(19, 7).Value = -17.33501518566552
$ws.Cells.Item(19, 8).Value = 1.555881390792372
$ws.Cells.Item(19, 9).Value = 1.694490435167829
$ws.Cells.Item(19, 10).Value = -17.33501518566552
$ws.Cells.Item(19, 11).Value = -17.33501518566552
$ws.Cells.Item(20, 2).Value = -17.33501518566552
$ws.Cells.Item(20, 3).Value = -17.33501518566552
$ws.Cells.Item(20, 4).Value = 2.017797691872143
$ws.Cells.Item(20, 5).Value = -17.33501518566552
$ws.Cells.Item(20, 6).Value = 3.221450641406618
$ws.Cells.Item(20, 7).Value = -17.33501518566552
$ws.Cells.Item(20, 8).Value = 2.065983041978745
$ws.Cells.Item(20, 9).Value = 1.021630210748218
$ws.Cells.Item(20, 10).Value = -17.33501518566552
$ws.Cells.Item(20, 11).Value = 2.295858031235161
$ws.Cells.Item(21, 2).Value = -17.33501518566552
$ws.Cells.Item(21, 3).Value = -17.33501518566552
$ws.Cells.Item(21, 4).Value = -17.33501518566552
$ws.Cells.Item(21, 5).Value = 2.553741438631586
$ws.Cells.Item(21, 6).Value = -17.33501518566552
$ws.Cells.Item(21, 7).Value = 2.550438410645449
$ws.Cells.Item(21, 8).Value = 2.27642933252156
$ws.Cells.Item(21, 9).Value = -17.33501518566552
$ws.Cells.Item(21, 10).Value = -17.33501518566552
$ws.Cells.Item(21, 11).Value = -17.33501518566552
